# Update the Rspo3-Lgr6 LR-pairs sheet with new TPM-derived values.
# Adds the "ECs" sending/target cluster alongside the existing "FAPs"
# cluster, giving the full FAPs/ECs x FAPs/ECs/MuSCs cross for the
# Rspo3 -> Lgr6 ligand-receptor pair (rows 2-7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column order: Sending cluster, Ligand symbol, Receptor symbol, Target
# cluster, then the 16 numeric NATMI metric columns (E..T).
$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$rowsData = @{
    2 = @("FAPs","Rspo3","Lgr6","FAPs",
          1, 0.3333333333333333, 0.01818866666666667, 0.054566,
          0.006403810693375696, 0.006403810693375696,
          1, 0.3333333333333333, 0.01182833333333333, 0.035485,
          0.03045376408867423, 0.03045376408867423,
          0.0002151416122222223, 0.00193627451,
          0.0001950201401245927, 0.0001950201401245927)
    3 = @("FAPs","Rspo3","Lgr6","ECs",
          1, 0.3333333333333333, 0.01818866666666667, 0.054566,
          0.006403810693375696, 0.006403810693375696,
          3, 1, 0.2280253333333333, 0.684076,
          0.5870843771374921, 0.5870843771374921,
          0.004147476779555556, 0.03732729101600001,
          0.003759577212226882, 0.003759577212226882)
    4 = @("FAPs","Rspo3","Lgr6","MuSCs",
          1, 0.3333333333333333, 0.01818866666666667, 0.054566,
          0.006403810693375696, 0.006403810693375696,
          3, 1, 0.1485493333333333, 0.445648,
          0.3824618587738337, 0.3824618587738337,
          0.002701914307555556, 0.024317228768,
          0.002449213341024221, 0.002449213341024221)
    5 = @("ECs","Rspo3","Lgr6","FAPs",
          3, 1, 2.822099333333334, 8.466298,
          0.9935961893066243, 0.9935961893066244,
          1, 0.3333333333333333, 0.01182833333333333, 0.035485,
          0.03045376408867423, 0.03045376408867423,
          0.03338073161444445, 0.30042658453,
          0.03025874394854963, 0.03025874394854964)
    6 = @("ECs","Rspo3","Lgr6","ECs",
          3, 1, 2.822099333333334, 8.466298,
          0.9935961893066243, 0.9935961893066244,
          3, 1, 0.2280253333333333, 0.684076,
          0.5870843771374921, 0.5870843771374921,
          0.6435101411831111, 5.791591270648,
          0.5833247999252652, 0.5833247999252652)
    7 = @("ECs","Rspo3","Lgr6","MuSCs",
          3, 1, 2.822099333333334, 8.466298,
          0.9935961893066243, 0.9935961893066244,
          3, 1, 0.1485493333333333, 0.445648,
          0.3824618587738337, 0.3824618587738337,
          0.4192209745671112, 3.772988771104,
          0.3800126454328095, 0.3800126454328095)
}

foreach ($r in 2..7) {
    $values = $rowsData[$r]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range("$($columns[$i])$r").Value = $values[$i]
    }
}
